# Move the picture on slide 3 ("Rho vs. Mu") up slightly.
# Old top = 1825625 EMU (143.75 pt); new top = 1690688 EMU (133.125... pt).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)

$shp.Top = 1690688 / 12700
